$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1, matching the style of the existing header row (G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data column values
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
